$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: unlock columns A:C entirely (column-level default + all existing cells)
$ws.Range("A1:C1000").Locked = $false

# Step 2: re-lock column A (labels) - all of column A
$ws.Range("A1:A1000").Locked = $true

# Step 3: re-lock the instructional note cell C3
$ws.Range("C3").Locked = $true

# Step 4: re-lock the border cell of the merged title row (B44)
$ws.Range("B44").Locked = $true

# Column B width change
$ws.Columns("B").ColumnWidth = 18.140625

# Protect the worksheet
$ws.Protect("", $true, $true, $true)

# Move selection to A44 (merged A44:B44) and scroll to top
$ws.Range("A44").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
